$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking Price/Volume values so they
# remain text cells (matching the source data) instead of being
# auto-converted to Number/Percentage by Excel's type inference.
$textCells = @("D2","D3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","E26","D38","E38","D39","E39","D40","E40","E41","D42","E42","D43","E43","D44","E44","D45","E45","E46","D47","E47","D48","D49","E49","D50","E50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell.
$ws.Range("D2").Value = '329.02'
$ws.Range("D3").Value = '44.41'
$ws.Range("D4").Value = '5.603'
$ws.Range("E4").Value = '3.58%'
$ws.Range("D5").Value = '0.08095'
$ws.Range("E5").Value = '0.07%'
$ws.Range("D6").Value = '2.039'
$ws.Range("E6").Value = '6.59%'
$ws.Range("D7").Value = '4.306'
$ws.Range("E7").Value = '0.05%'
$ws.Range("D8").Value = '0.9531'
$ws.Range("E8").Value = '1.15%'
$ws.Range("D9").Value = '2.564'
$ws.Range("E9").Value = '-7.29%'
$ws.Range("D10").Value = '0.1174'
$ws.Range("E10").Value = '-2.77%'
$ws.Range("D11").Value = '0.1855'
$ws.Range("E11").Value = '-2.40%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '0.09785'
$ws.Range("E12").Value = '2.78%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.04579'
$ws.Range("E13").Value = '9.35%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.1068'
$ws.Range("E14").Value = '-0.07%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001283'
$ws.Range("E15").Value = '-2.00%'
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").Value = '0.04210'
$ws.Range("E16").Value = '-3.83%'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = '0.005856'
$ws.Range("E17").Value = '-2.69%'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = '3.368'
$ws.Range("E18").Value = '-5.64%'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '0.3474'
$ws.Range("E19").Value = '-0.71%'
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D20").Value = '10.18'
$ws.Range("E20").Value = '19.54%'
$ws.Range("D21").Value = '0.1408'
$ws.Range("E21").Value = '3.71%'
$ws.Range("D22").Value = '0.2503'
$ws.Range("E22").Value = '-3.86%'
$ws.Range("D23").Value = '0.001245'
$ws.Range("E23").Value = '0.42%'
$ws.Range("D24").Value = '0.004318'
$ws.Range("E24").Value = '0.18%'
$ws.Range("D25").Value = '0.0001189'
$ws.Range("E25").Value = '-4.00%'
$ws.Range("E26").Value = '-0.88%'
$ws.Range("D38").Value = '0.02669'
$ws.Range("E38").Value = '-0.06%'
$ws.Range("D39").Value = '0.05556'
$ws.Range("E39").Value = '1.80%'
$ws.Range("D40").Value = '0.007582'
$ws.Range("E40").Value = '-2.79%'
$ws.Range("E41").Value = '1.07%'
$ws.Range("D42").Value = '0.008072'
$ws.Range("E42").Value = '-17.38%'
$ws.Range("D43").Value = '0.002014'
$ws.Range("E43").Value = '-5.27%'
$ws.Range("D44").Value = '0.008384'
$ws.Range("E44").Value = '-15.87%'
$ws.Range("D45").Value = '0.00007204'
$ws.Range("E45").Value = '1.82%'
$ws.Range("E46").Value = '-0.72%'
$ws.Range("D47").Value = '0.004156'
$ws.Range("E47").Value = '19.68%'
$ws.Range("D48").Value = '0.002269'
$ws.Range("D49").Value = '0.00002099'
$ws.Range("E49").Value = '-0.72%'
$ws.Range("D50").Value = '0.0001999'
$ws.Range("E50").Value = '-0.72%'
